$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Rename the product name (shared by both sheets' B1 cell) to remove the
# test-case inter-dependency on the literal "...PERIODIC" suffix.
$newProductName = "2575-MS-EPP-DB-SAR-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-PER-1st"
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# Shortname on the input sheet becomes a literal string instead of the
# numeric 2575 value.
$ws1.Range("B2").Value = "257e"

# Reset selection on the input sheet back to B1 and make it no longer the
# active/selected tab.
$ws1.Range("B1").Select()

# Make the output sheet the active sheet/tab instead.
$ws2.Activate()
$ws2.Range("B1").Select()
